# Auto-generated edit script: updates Kraken market-data values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 49
$ws.Range("I6").Value = 38.8
$ws.Range("K6").Value = 116.4
$ws.Range("M6").Value = -4.399999999999991
$ws.Range("H9").Value = 356.125
$ws.Range("I9").Value = 356.125
$ws.Range("K9").Value = 356.125
$ws.Range("M9").Value = -187.125
$ws.Range("H33").Value = 235.55556
$ws.Range("I33").Value = 235.55556
$ws.Range("K33").Value = 235.55556
$ws.Range("M33").Value = -6.555560000000014
$ws.Range("H98").Value = 3777.5
$ws.Range("I98").Value = 4181
$ws.Range("J98").Value = 1760
$ws.Range("K98").Value = 4181
$ws.Range("L98").Value = 1760
$ws.Range("M98").Value = -2683
$ws.Range("N98").Value = -4756
$ws.Range("H107").Value = 2068.25
$ws.Range("I107").Value = 1341.1666
$ws.Range("K107").Value = 1341.1666
$ws.Range("M107").Value = 578.8334
$ws.Range("H122").Value = 3777.5
$ws.Range("I122").Value = 4181
$ws.Range("J122").Value = 1760
$ws.Range("K122").Value = 12543
$ws.Range("L122").Value = 5280
$ws.Range("M122").Value = -10093
$ws.Range("N122").Value = -10180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1049.5
$ws.Range("I45").Value = 999.5
$ws.Range("J45").Value = 1099.5
$ws.Range("K45").Value = 999.5
$ws.Range("L45").Value = 1099.5
$ws.Range("M45").Value = -622.5
$ws.Range("N45").Value = -1853.5
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("H110").Value = 1102.5
$ws.Range("I110").Value = 705.5
$ws.Range("J110").Value = 1499.5
$ws.Range("K110").Value = 705.5
$ws.Range("L110").Value = 1499.5
$ws.Range("M110").Value = 1339.5
$ws.Range("N110").Value = -5589.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1251.625
$ws.Range("I86").Value = 1143.9166
$ws.Range("K86").Value = 1143.9166
$ws.Range("M86").Value = -20.91660000000002
$ws.Range("H89").Value = 1251.625
$ws.Range("I89").Value = 1143.9166
$ws.Range("K89").Value = 5719.583000000001
$ws.Range("M89").Value = -103.5830000000005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 255.11111
$ws.Range("I7").Value = 179.6
$ws.Range("J7").Value = 349.5
$ws.Range("K7").Value = 179.6
$ws.Range("L7").Value = 349.5
$ws.Range("M7").Value = -66.59999999999999
$ws.Range("N7").Value = -575.5
$ws.Range("H31").Value = 2758.875
$ws.Range("I31").Value = 1012
$ws.Range("K31").Value = 1012
$ws.Range("M31").Value = -717
$ws.Range("H34").Value = 2758.875
$ws.Range("I34").Value = 1012
$ws.Range("K34").Value = 1012
$ws.Range("M34").Value = -810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 799999
$ws.Range("I58").Value = 799999
$ws.Range("K58").Value = 799999
$ws.Range("M58").Value = -799722
$ws.Range("H80").Value = 5133.3335
$ws.Range("J80").Value = 4950
$ws.Range("L80").Value = 4950
$ws.Range("N80").Value = -6946
$ws.Range("H83").Value = 5133.3335
$ws.Range("J83").Value = 4950
$ws.Range("L83").Value = 24750
$ws.Range("N83").Value = -34734
$ws.Range("H113").Value = 5959.8
$ws.Range("I113").Value = 4900
$ws.Range("J113").Value = 6666.3335
$ws.Range("K113").Value = 4900
$ws.Range("L113").Value = 6666.3335
$ws.Range("M113").Value = -2730
$ws.Range("N113").Value = -11006.3335
$ws.Range("H122").Value = 1211.25
$ws.Range("I122").Value = 1165
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 3495
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -1045
$ws.Range("N122").Value = -8950
$ws.Range("H126").Value = 3942.1428
$ws.Range("I126").Value = 4074.1667
$ws.Range("J126").Value = 3150
$ws.Range("K126").Value = 12222.5001
$ws.Range("L126").Value = 9450
$ws.Range("M126").Value = -9752.500100000001
$ws.Range("N126").Value = -14390
$ws.Range("H132").Value = 2459.3635
$ws.Range("I132").Value = 2176.5715
$ws.Range("J132").Value = 2954.25
$ws.Range("K132").Value = 6529.7145
$ws.Range("L132").Value = 8862.75
$ws.Range("M132").Value = -3999.7145
$ws.Range("N132").Value = -13922.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4160.5557
$ws.Range("J7").Value = 4749.25
$ws.Range("L7").Value = 4749.25
$ws.Range("N7").Value = -4973.25
$ws.Range("H61").Value = 3816.3333
$ws.Range("I61").Value = 3474.75
$ws.Range("K61").Value = 3474.75
$ws.Range("M61").Value = -3272.75
$ws.Range("H82").Value = 1800
$ws.Range("I82").Value = 1650
$ws.Range("K82").Value = 1650
$ws.Range("M82").Value = -1289
$ws.Range("H85").Value = 1800
$ws.Range("I85").Value = 1650
$ws.Range("K85").Value = 1650
$ws.Range("M85").Value = -402
$ws.Range("H100").Value = 8000
$ws.Range("I100").Value = 6000
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -5459
$ws.Range("N100").Value = -11082
$ws.Range("H113").Value = 3816.3333
$ws.Range("I113").Value = 3474.75
$ws.Range("K113").Value = 3474.75
$ws.Range("M113").Value = -1304.75
$ws.Range("H126").Value = 4160.5557
$ws.Range("J126").Value = 4749.25
$ws.Range("L126").Value = 14247.75
$ws.Range("N126").Value = -19187.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000000
$ws.Range("J5").Value = 10000000
$ws.Range("L5").Value = 10000000
$ws.Range("N5").Value = -10000224
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H126").Value = 1888.8334
$ws.Range("I126").Value = 1888.8334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5666.5002
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -3196.5002
